$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for Price/Volume (and swapped Coin/Link) cells so numeric-looking
# strings like "1.002" are not reinterpreted as numbers.
$ws.Range("D2:E51").NumberFormat = "@"
$ws.Range("B33:C34").NumberFormat = "@"
$ws.Range("B38:C39").NumberFormat = "@"

$ws.Range("D2").Value = "22.426.38"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").Value = "1.568.51"
$ws.Range("E3").Value = "  -0.30%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "1.002"
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("D6").Value = "285.27"
$ws.Range("E6").Value = "  -2.23%  "
$ws.Range("D7").Value = "0.3632"
$ws.Range("E7").Value = "  -2.53%  "
$ws.Range("D8").Value = "48.58"
$ws.Range("E8").Value = "  -2.74%  "
$ws.Range("D9").Value = "0.3328"
$ws.Range("E9").Value = "  -1.90%  "
$ws.Range("D10").Value = "1.124"
$ws.Range("E10").Value = "  -1.86%  "
$ws.Range("D11").Value = "0.07393"
$ws.Range("E11").Value = "  -2.36%  "
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("D13").Value = "20.78"
$ws.Range("E13").Value = "  -2.34%  "
$ws.Range("D14").Value = "5.952"
$ws.Range("E14").Value = "  -1.06%  "
$ws.Range("D15").Value = "6.907"
$ws.Range("E15").Value = "  -0.74%  "
$ws.Range("D16").Value = "1.569.70"
$ws.Range("E16").Value = "  -0.25%  "
$ws.Range("D17").Value = "0.00001104"
$ws.Range("E17").Value = "  -1.65%  "
$ws.Range("D18").Value = "88.06"
$ws.Range("E18").Value = "  -3.21%  "
$ws.Range("D19").Value = "0.06698"
$ws.Range("E19").Value = "  -0.84%  "
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").Value = "6.340"
$ws.Range("E21").Value = "  +0.61%  "
$ws.Range("D22").Value = "16.19"
$ws.Range("E22").Value = "  -0.74%  "
$ws.Range("D23").Value = "12.01"
$ws.Range("E23").Value = "  -1.26%  "
$ws.Range("D24").Value = "22.414.78"
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("D25").Value = "2.382"
$ws.Range("E25").Value = "  +1.86%  "
$ws.Range("D26").Value = "2.539"
$ws.Range("E26").Value = "  -5.65%  "
$ws.Range("D27").Value = "150.55"
$ws.Range("E27").Value = "  +1.22%  "
$ws.Range("D28").Value = "19.41"
$ws.Range("E28").Value = "  -3.34%  "
$ws.Range("D29").Value = "5.002"
$ws.Range("E29").Value = "  -0.39%  "
$ws.Range("D30").Value = "123.86"
$ws.Range("E30").Value = "  -1.45%  "
$ws.Range("D31").Value = "1.745.97"
$ws.Range("E31").Value = "  -0.30%  "
$ws.Range("D32").Value = "1.042"
$ws.Range("E32").Value = "  -1.98%  "
$ws.Range("D33").Value = "2.000"
$ws.Range("E33").Value = "  +0.65%  "
$ws.Range("D34").Value = "6.099"
$ws.Range("E34").Value = "  -1.16%  "
$ws.Range("D35").Value = "9.803"
$ws.Range("E35").Value = "  -0.45%  "
$ws.Range("D36").Value = "0.08243"
$ws.Range("E36").Value = "  -1.49%  "
$ws.Range("D37").Value = "0.02414"
$ws.Range("E37").Value = "  -3.18%  "
$ws.Range("D38").Value = "0.2234"
$ws.Range("E38").Value = "  -3.09%  "
$ws.Range("D39").Value = "0.06426"
$ws.Range("E39").Value = "  -1.42%  "
$ws.Range("D40").Value = "5.359"
$ws.Range("E40").Value = "  -1.96%  "
$ws.Range("D41").Value = "1.287"
$ws.Range("E41").Value = "  -5.14%  "
$ws.Range("D42").Value = "0.6276"
$ws.Range("E42").Value = "  +0.95%  "
$ws.Range("D43").Value = "11.19"
$ws.Range("E43").Value = "  -1.12%  "
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").Value = "13.88"
$ws.Range("E45").Value = "  -0.77%  "
$ws.Range("D46").Value = "0.6074"
$ws.Range("E46").Value = "  +4.63%  "
$ws.Range("D47").Value = "3.751"
$ws.Range("E47").Value = "  -1.68%  "
$ws.Range("D48").Value = "2.034"
$ws.Range("E48").Value = "  -1.55%  "
$ws.Range("D49").Value = "123.51"
$ws.Range("E49").Value = "  -5.00%  "
$ws.Range("D50").Value = "1.212"
$ws.Range("E50").Value = "  -0.93%  "
$ws.Range("D51").Value = "0.07209"
$ws.Range("E51").Value = "  -1.56%  "

# Rows 33/34 and 38/39: coin identity (Coin name + Link) swapped between the two rows
$ws.Range("B33").Value = "WEMIXTOKEN"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"

# Restore default (unformatted) style so only cell content changed, matching the source
# workbook which never applies an explicit number format to these cells.
$ws.Range("D2:E51").Style = "Normal"
$ws.Range("B33:C34").Style = "Normal"
$ws.Range("B38:C39").Style = "Normal"
